$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.582
$ws.Range("E4").Value = 13.265

$ws.Range("E5").Value = 13.254

$ws.Range("A7").Value = -20.987

$ws.Range("E8").Value = 13.718

$ws.Range("A16").Value = -20.806
$ws.Range("E16").Value = 13.146
